$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "Save" column - copy formatting from the adjacent "sum" header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save values per row (row 1 is header, data in rows 2-14)
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
